$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("hall_id") ahead of the existing stage_size /
# no_of_green_rooms columns, shifting the old data right.
$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "hall_id"

# hall_id values were typed with a leading apostrophe (text-quote prefix)
# plus literal leading/trailing single-quote characters, e.g. 'H03' -
# doubling the leading apostrophe keeps one as the literal character while
# the first still acts as Excel's "force text" prefix.
$ws.Range("A2").Value = "''H03'"
$ws.Range("A3").Value = "''H04'"
$ws.Range("A4").Value = "''H06'"
$ws.Range("A5").Value = "''H08'"
$ws.Range("A6").Value = "''H10'"
$ws.Range("A7").Value = "''H11'"
$ws.Range("A8").Value = "''H14'"
$ws.Range("A9").Value = "''H15'"
$ws.Range("A10").Value = "''H19'"
$ws.Range("A11").Value = "''H20'"

$ws.Range("B2").Value = "''2000 sq feet'"
$ws.Range("B3").Value = "''1500 sq feet'"
$ws.Range("B4").Value = "''2000 sq feet'"
$ws.Range("B5").Value = "''3000 sq feet'"
$ws.Range("B6").Value = "''3000 sq feet'"
$ws.Range("B7").Value = "''2000 sq feet'"
$ws.Range("B8").Value = "''3500 sq feet'"
$ws.Range("B9").Value = "''3500 sq feet'"
$ws.Range("B10").Value = "''2000 sq feet'"
$ws.Range("B11").Value = "''1500 sq feet'"

$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 3
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 2

# Match the selection cell the author's session ended on.
$null = $ws.Range("F12").Select()
